$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new ingredient row for "curry powder" (egg curry update)
$newRow = 84

$ws.Cells.Item($newRow, 1).Value = "curry powder"
$ws.Cells.Item($newRow, 2).Value = "Check"
$ws.Cells.Item($newRow, 3).Value = 0
$ws.Cells.Item($newRow, 4).Value = 0
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0

# Update the view to reflect the newly selected/active cell and scroll position
$ws.Range("A84").Select()
$excel.ActiveWindow.ScrollRow = 64
